$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates per diff. D-column numeric-looking values are forced to Text
# format before assignment (then format is reset to Normal) so Excel does not
# silently convert them into numeric cells - matching the original inlineStr layout.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.838.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.815.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9976"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.44%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "337.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9917"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3937"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.74%  "
$ws.Range("E8").Value = "  +1.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.91"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07576"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9913"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.539"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.810.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.203"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001107"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06670"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "85.33"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9920"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.599"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.878.74"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.417"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.570"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.501"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "155.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.017.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "135.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.021"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.117"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08814"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.29"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.545"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.25%  "
$ws.Range("B37").Value = "WEMIXTOKEN"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.625"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6916"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.41%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02421"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.93%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06530"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2231"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.267"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.558"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.70"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6548"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9908"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.861"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.166"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07222"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.89"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.47%  "
